# Xpath scraping is done: the scraper appends newly scraped "FireFlink_*"
# project names to the workbook's shared-string history and repoints the
# "latest value" cell (Sheet1!B2) at the most recently scraped name.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Re-assert every previously scraped name (including the one Sheet1!B2 used
# to point at) on the scratch sheet so none of the accumulated history is
# lost, then add this run's 128 newly scraped names after them, in order.
$ws2.Range("A1").Value = "FireFlink_63430"
$ws2.Range("B1").Value = "FireFlink_14020"
$ws2.Range("C1").Value = "FireFlink_24392"
$ws2.Range("D1").Value = "FireFlink_47600"
$ws2.Range("E1").Value = "FireFlink_66868"
$ws2.Range("F1").Value = "FireFlink_86351"
$ws2.Range("G1").Value = "FireFlink_55296"
$ws2.Range("H1").Value = "FireFlink_73882"
$ws2.Range("I1").Value = "FireFlink_41275"
$ws2.Range("J1").Value = "FireFlink_04851"
$ws2.Range("K1").Value = "FireFlink_73498"
$ws2.Range("L1").Value = "FireFlink_60480"
$ws2.Range("M1").Value = "FireFlink_74969"
$ws2.Range("N1").Value = "FireFlink_27218"
$ws2.Range("O1").Value = "FireFlink_52885"
$ws2.Range("P1").Value = "FireFlink_88164"
$ws2.Range("Q1").Value = "FireFlink_96339"
$ws2.Range("R1").Value = "FireFlink_15707"
$ws2.Range("S1").Value = "FireFlink_00609"
$ws2.Range("T1").Value = "FireFlink_34359"
$ws2.Range("A2").Value = "FireFlink_02483"
$ws2.Range("B2").Value = "FireFlink_30557"
$ws2.Range("C2").Value = "FireFlink_94246"
$ws2.Range("D2").Value = "FireFlink_88894"
$ws2.Range("E2").Value = "FireFlink_01788"
$ws2.Range("F2").Value = "FireFlink_57975"
$ws2.Range("G2").Value = "FireFlink_15162"
$ws2.Range("H2").Value = "FireFlink_22749"
$ws2.Range("I2").Value = "FireFlink_71104"
$ws2.Range("J2").Value = "FireFlink_59641"
$ws2.Range("K2").Value = "FireFlink_88465"
$ws2.Range("L2").Value = "FireFlink_78770"
$ws2.Range("M2").Value = "FireFlink_12501"
$ws2.Range("N2").Value = "FireFlink_74229"
$ws2.Range("O2").Value = "FireFlink_82246"
$ws2.Range("P2").Value = "FireFlink_42796"
$ws2.Range("Q2").Value = "FireFlink_51421"
$ws2.Range("R2").Value = "FireFlink_85215"
$ws2.Range("S2").Value = "FireFlink_28403"
$ws2.Range("T2").Value = "FireFlink_54696"
$ws2.Range("A3").Value = "FireFlink_49324"
$ws2.Range("B3").Value = "FireFlink_20814"
$ws2.Range("C3").Value = "FireFlink_29820"
$ws2.Range("D3").Value = "FireFlink_14697"
$ws2.Range("E3").Value = "FireFlink_18301"
$ws2.Range("F3").Value = "FireFlink_79065"
$ws2.Range("G3").Value = "FireFlink_82094"
$ws2.Range("H3").Value = "FireFlink_99431"
$ws2.Range("I3").Value = "FireFlink_72675"
$ws2.Range("J3").Value = "FireFlink_51379"
$ws2.Range("K3").Value = "FireFlink_50115"
$ws2.Range("L3").Value = "FireFlink_80814"
$ws2.Range("M3").Value = "FireFlink_62358"
$ws2.Range("N3").Value = "FireFlink_08433"
$ws2.Range("O3").Value = "FireFlink_59371"
$ws2.Range("P3").Value = "FireFlink_36392"
$ws2.Range("Q3").Value = "FireFlink_25030"
$ws2.Range("R3").Value = "FireFlink_45093"
$ws2.Range("S3").Value = "FireFlink_07183"
$ws2.Range("T3").Value = "FireFlink_91060"
$ws2.Range("A4").Value = "FireFlink_44625"
$ws2.Range("B4").Value = "FireFlink_94642"
$ws2.Range("C4").Value = "FireFlink_84988"
$ws2.Range("D4").Value = "FireFlink_40555"
$ws2.Range("E4").Value = "FireFlink_25908"
$ws2.Range("F4").Value = "FireFlink_31501"
$ws2.Range("G4").Value = "FireFlink_13681"
$ws2.Range("H4").Value = "FireFlink_73133"
$ws2.Range("I4").Value = "FireFlink_43456"
$ws2.Range("J4").Value = "FireFlink_95456"
$ws2.Range("K4").Value = "FireFlink_20862"
$ws2.Range("L4").Value = "FireFlink_79988"
$ws2.Range("M4").Value = "FireFlink_94557"
$ws2.Range("N4").Value = "FireFlink_12102"
$ws2.Range("O4").Value = "FireFlink_73640"
$ws2.Range("P4").Value = "FireFlink_62538"
$ws2.Range("Q4").Value = "FireFlink_45755"
$ws2.Range("R4").Value = "FireFlink_25537"
$ws2.Range("S4").Value = "FireFlink_84685"
$ws2.Range("T4").Value = "FireFlink_49073"
$ws2.Range("A5").Value = "FireFlink_56027"
$ws2.Range("B5").Value = "FireFlink_91709"
$ws2.Range("C5").Value = "FireFlink_89788"
$ws2.Range("D5").Value = "FireFlink_25605"
$ws2.Range("E5").Value = "FireFlink_03619"
$ws2.Range("F5").Value = "FireFlink_75137"
$ws2.Range("G5").Value = "FireFlink_79034"
$ws2.Range("H5").Value = "FireFlink_55608"
$ws2.Range("I5").Value = "FireFlink_45155"
$ws2.Range("J5").Value = "FireFlink_31839"
$ws2.Range("K5").Value = "FireFlink_84367"
$ws2.Range("L5").Value = "FireFlink_63020"
$ws2.Range("M5").Value = "FireFlink_73231"
$ws2.Range("N5").Value = "FireFlink_27330"
$ws2.Range("O5").Value = "FireFlink_65983"
$ws2.Range("P5").Value = "FireFlink_81727"
$ws2.Range("Q5").Value = "FireFlink_44382"
$ws2.Range("R5").Value = "FireFlink_49011"
$ws2.Range("S5").Value = "FireFlink_41826"
$ws2.Range("T5").Value = "FireFlink_83074"
$ws2.Range("A6").Value = "FireFlink_55484"
$ws2.Range("B6").Value = "FireFlink_15572"
$ws2.Range("C6").Value = "FireFlink_76640"
$ws2.Range("D6").Value = "FireFlink_20757"
$ws2.Range("E6").Value = "FireFlink_25309"
$ws2.Range("F6").Value = "FireFlink_27457"
$ws2.Range("G6").Value = "FireFlink_12642"
$ws2.Range("H6").Value = "FireFlink_42880"
$ws2.Range("I6").Value = "FireFlink_68878"
$ws2.Range("J6").Value = "FireFlink_54388"
$ws2.Range("K6").Value = "FireFlink_82827"
$ws2.Range("L6").Value = "FireFlink_80666"
$ws2.Range("M6").Value = "FireFlink_38778"
$ws2.Range("N6").Value = "FireFlink_24433"
$ws2.Range("O6").Value = "FireFlink_34146"
$ws2.Range("P6").Value = "FireFlink_43122"
$ws2.Range("Q6").Value = "FireFlink_69018"
$ws2.Range("R6").Value = "FireFlink_50752"
$ws2.Range("S6").Value = "FireFlink_04719"
$ws2.Range("T6").Value = "FireFlink_48464"
$ws2.Range("A7").Value = "FireFlink_88654"
$ws2.Range("B7").Value = "FireFlink_24426"
$ws2.Range("C7").Value = "FireFlink_64446"
$ws2.Range("D7").Value = "FireFlink_04277"
$ws2.Range("E7").Value = "FireFlink_42919"
$ws2.Range("F7").Value = "FireFlink_68020"
$ws2.Range("G7").Value = "FireFlink_10418"
$ws2.Range("H7").Value = "FireFlink_31278"
$ws2.Range("I7").Value = "FireFlink_09012"
$ws2.Range("J7").Value = "FireFlink_08897"
$ws2.Range("K7").Value = "FireFlink_74597"
$ws2.Range("L7").Value = "FireFlink_27813"
$ws2.Range("M7").Value = "FireFlink_57140"
$ws2.Range("N7").Value = "FireFlink_01594"
$ws2.Range("O7").Value = "FireFlink_42877"
$ws2.Range("P7").Value = "FireFlink_26368"
$ws2.Range("Q7").Value = "FireFlink_36866"
$ws2.Range("R7").Value = "FireFlink_51665"
$ws2.Range("S7").Value = "FireFlink_69365"
$ws2.Range("T7").Value = "FireFlink_58899"
$ws2.Range("A8").Value = "FireFlink_26459"
$ws2.Range("B8").Value = "FireFlink_72540"
$ws2.Range("C8").Value = "FireFlink_06792"
$ws2.Range("D8").Value = "FireFlink_46702"
$ws2.Range("E8").Value = "FireFlink_50126"
$ws2.Range("F8").Value = "FireFlink_96862"
$ws2.Range("G8").Value = "FireFlink_48714"
$ws2.Range("H8").Value = "FireFlink_63756"
$ws2.Range("I8").Value = "FireFlink_55354"
$ws2.Range("J8").Value = "FireFlink_82692"
$ws2.Range("K8").Value = "FireFlink_11492"
$ws2.Range("L8").Value = "FireFlink_11882"
$ws2.Range("M8").Value = "FireFlink_99041"
$ws2.Range("N8").Value = "FireFlink_39571"
$ws2.Range("O8").Value = "FireFlink_58304"
$ws2.Range("P8").Value = "FireFlink_79605"
$ws2.Range("Q8").Value = "FireFlink_73107"
$ws2.Range("R8").Value = "FireFlink_79075"
$ws2.Range("S8").Value = "FireFlink_83378"
$ws2.Range("T8").Value = "FireFlink_27644"
$ws2.Range("A9").Value = "FireFlink_92580"
$ws2.Range("B9").Value = "FireFlink_14300"
$ws2.Range("C9").Value = "FireFlink_05073"
$ws2.Range("D9").Value = "FireFlink_19390"
$ws2.Range("E9").Value = "FireFlink_35143"
$ws2.Range("F9").Value = "FireFlink_62522"
$ws2.Range("G9").Value = "FireFlink_74524"
$ws2.Range("H9").Value = "FireFlink_05846"
$ws2.Range("I9").Value = "FireFlink_24938"
$ws2.Range("J9").Value = "FireFlink_31255"
$ws2.Range("K9").Value = "FireFlink_02690"
$ws2.Range("L9").Value = "FireFlink_58335"
$ws2.Range("M9").Value = "FireFlink_25755"
$ws2.Range("N9").Value = "FireFlink_35949"
$ws2.Range("O9").Value = "FireFlink_37604"
$ws2.Range("P9").Value = "FireFlink_17713"
$ws2.Range("Q9").Value = "FireFlink_22283"
$ws2.Range("R9").Value = "FireFlink_37139"
$ws2.Range("S9").Value = "FireFlink_60242"
$ws2.Range("T9").Value = "FireFlink_52907"
$ws2.Range("A10").Value = "FireFlink_50315"
$ws2.Range("B10").Value = "FireFlink_83401"
$ws2.Range("C10").Value = "FireFlink_28489"
$ws2.Range("D10").Value = "FireFlink_69113"
$ws2.Range("E10").Value = "FireFlink_51576"
$ws2.Range("F10").Value = "FireFlink_39640"
$ws2.Range("G10").Value = "FireFlink_18099"
$ws2.Range("H10").Value = "FireFlink_31582"
$ws2.Range("I10").Value = "FireFlink_08127"
$ws2.Range("J10").Value = "FireFlink_19976"
$ws2.Range("K10").Value = "FireFlink_80395"
$ws2.Range("L10").Value = "FireFlink_10990"
$ws2.Range("M10").Value = "FireFlink_94103"
$ws2.Range("N10").Value = "FireFlink_55412"
$ws2.Range("O10").Value = "FireFlink_25520"
$ws2.Range("P10").Value = "FireFlink_25438"
$ws2.Range("Q10").Value = "FireFlink_35339"
$ws2.Range("R10").Value = "FireFlink_86900"
$ws2.Range("S10").Value = "FireFlink_77357"
$ws2.Range("T10").Value = "FireFlink_24502"
$ws2.Range("A11").Value = "FireFlink_18027"
$ws2.Range("B11").Value = "FireFlink_96267"
$ws2.Range("C11").Value = "FireFlink_42725"
$ws2.Range("D11").Value = "FireFlink_67087"
$ws2.Range("E11").Value = "FireFlink_65721"
$ws2.Range("F11").Value = "FireFlink_11669"
$ws2.Range("G11").Value = "FireFlink_23558"
$ws2.Range("H11").Value = "FireFlink_07337"
$ws2.Range("I11").Value = "FireFlink_35835"
$ws2.Range("J11").Value = "FireFlink_66807"
$ws2.Range("K11").Value = "FireFlink_58295"
$ws2.Range("L11").Value = "FireFlink_90147"
$ws2.Range("M11").Value = "FireFlink_02512"
$ws2.Range("N11").Value = "FireFlink_62325"
$ws2.Range("O11").Value = "FireFlink_56485"
$ws2.Range("P11").Value = "FireFlink_20187"
$ws2.Range("Q11").Value = "FireFlink_76439"
$ws2.Range("R11").Value = "FireFlink_44962"
$ws2.Range("S11").Value = "FireFlink_88147"
$ws2.Range("T11").Value = "FireFlink_30438"
$ws2.Range("A12").Value = "FireFlink_21760"
$ws2.Range("B12").Value = "FireFlink_12015"
$ws2.Range("C12").Value = "FireFlink_02681"
$ws2.Range("D12").Value = "FireFlink_93864"
$ws2.Range("E12").Value = "FireFlink_69080"
$ws2.Range("F12").Value = "FireFlink_63678"
$ws2.Range("G12").Value = "FireFlink_97963"
$ws2.Range("H12").Value = "FireFlink_96888"
$ws2.Range("I12").Value = "FireFlink_53429"
$ws2.Range("J12").Value = "FireFlink_38738"
$ws2.Range("K12").Value = "FireFlink_14465"
$ws2.Range("L12").Value = "FireFlink_08039"
$ws2.Range("M12").Value = "FireFlink_17246"
$ws2.Range("N12").Value = "FireFlink_10848"
$ws2.Range("O12").Value = "FireFlink_58232"
$ws2.Range("P12").Value = "FireFlink_90040"
$ws2.Range("Q12").Value = "FireFlink_59825"
$ws2.Range("R12").Value = "FireFlink_02636"
$ws2.Range("S12").Value = "FireFlink_31289"
$ws2.Range("T12").Value = "FireFlink_58321"
$ws2.Range("A13").Value = "FireFlink_47846"
$ws2.Range("B13").Value = "FireFlink_68993"
$ws2.Range("C13").Value = "FireFlink_43568"
$ws2.Range("D13").Value = "FireFlink_43756"
$ws2.Range("E13").Value = "FireFlink_68775"
$ws2.Range("F13").Value = "FireFlink_21879"
$ws2.Range("G13").Value = "FireFlink_46836"
$ws2.Range("H13").Value = "FireFlink_62983"
$ws2.Range("I13").Value = "FireFlink_22755"
$ws2.Range("J13").Value = "FireFlink_95290"
$ws2.Range("K13").Value = "FireFlink_76655"
$ws2.Range("L13").Value = "FireFlink_23017"
$ws2.Range("M13").Value = "FireFlink_57305"
$ws2.Range("N13").Value = "FireFlink_94349"
$ws2.Range("O13").Value = "FireFlink_36314"
$ws2.Range("P13").Value = "FireFlink_64099"
$ws2.Range("Q13").Value = "FireFlink_05764"
$ws2.Range("R13").Value = "FireFlink_00445"
$ws2.Range("S13").Value = "FireFlink_26865"
$ws2.Range("T13").Value = "FireFlink_17481"
$ws2.Range("A14").Value = "FireFlink_53830"
$ws2.Range("B14").Value = "FireFlink_63199"
$ws2.Range("C14").Value = "FireFlink_90693"
$ws2.Range("D14").Value = "FireFlink_60086"
$ws2.Range("E14").Value = "FireFlink_50003"
$ws2.Range("F14").Value = "FireFlink_41409"
$ws2.Range("G14").Value = "FireFlink_56052"
$ws2.Range("H14").Value = "FireFlink_97288"
$ws2.Range("I14").Value = "FireFlink_09398"
$ws2.Range("J14").Value = "FireFlink_86570"
$ws2.Range("K14").Value = "FireFlink_95223"
$ws2.Range("L14").Value = "FireFlink_97363"
$ws2.Range("M14").Value = "FireFlink_28588"
$ws2.Range("N14").Value = "FireFlink_54457"
$ws2.Range("O14").Value = "FireFlink_12167"
$ws2.Range("P14").Value = "FireFlink_55236"
$ws2.Range("Q14").Value = "FireFlink_32663"
$ws2.Range("R14").Value = "FireFlink_64894"
$ws2.Range("S14").Value = "FireFlink_51760"
$ws2.Range("T14").Value = "FireFlink_59997"
$ws2.Range("A15").Value = "FireFlink_71894"
$ws2.Range("B15").Value = "FireFlink_68169"
$ws2.Range("C15").Value = "FireFlink_23447"
$ws2.Range("D15").Value = "FireFlink_31194"
$ws2.Range("E15").Value = "FireFlink_45698"
$ws2.Range("F15").Value = "FireFlink_22102"
$ws2.Range("G15").Value = "FireFlink_41559"
$ws2.Range("H15").Value = "FireFlink_92950"
$ws2.Range("I15").Value = "FireFlink_18004"
$ws2.Range("J15").Value = "FireFlink_48778"
$ws2.Range("K15").Value = "FireFlink_32063"
$ws2.Range("L15").Value = "FireFlink_51439"
$ws2.Range("M15").Value = "FireFlink_87345"
$ws2.Range("N15").Value = "FireFlink_41674"
$ws2.Range("O15").Value = "FireFlink_14109"
$ws2.Range("P15").Value = "FireFlink_60141"
$ws2.Range("Q15").Value = "FireFlink_19079"
$ws2.Range("R15").Value = "FireFlink_87032"
$ws2.Range("S15").Value = "FireFlink_74709"
$ws2.Range("T15").Value = "FireFlink_47030"
$ws2.Range("A16").Value = "FireFlink_69893"
$ws2.Range("B16").Value = "FireFlink_36833"
$ws2.Range("C16").Value = "FireFlink_58887"
$ws2.Range("D16").Value = "FireFlink_57224"
$ws2.Range("E16").Value = "FireFlink_86092"
$ws2.Range("F16").Value = "FireFlink_02231"
$ws2.Range("G16").Value = "FireFlink_62562"
$ws2.Range("H16").Value = "FireFlink_11863"
$ws2.Range("I16").Value = "FireFlink_87274"
$ws2.Range("J16").Value = "FireFlink_59202"
$ws2.Range("K16").Value = "FireFlink_11757"
$ws2.Range("L16").Value = "FireFlink_04964"
$ws2.Range("M16").Value = "FireFlink_48302"
$ws2.Range("N16").Value = "FireFlink_96633"
$ws2.Range("O16").Value = "FireFlink_60974"
$ws2.Range("P16").Value = "FireFlink_71085"
$ws2.Range("Q16").Value = "FireFlink_19109"
$ws2.Range("R16").Value = "FireFlink_30507"
$ws2.Range("S16").Value = "FireFlink_93894"
$ws2.Range("T16").Value = "FireFlink_10115"
$ws2.Range("A17").Value = "FireFlink_37206"
$ws2.Range("B17").Value = "FireFlink_91585"
$ws2.Range("C17").Value = "FireFlink_01932"
$ws2.Range("D17").Value = "FireFlink_93781"
$ws2.Range("E17").Value = "FireFlink_67971"
$ws2.Range("F17").Value = "FireFlink_98663"
$ws2.Range("G17").Value = "FireFlink_93558"
$ws2.Range("H17").Value = "FireFlink_99713"
$ws2.Range("I17").Value = "FireFlink_27908"
$ws2.Range("J17").Value = "FireFlink_88820"
$ws2.Range("K17").Value = "FireFlink_58646"
$ws2.Range("L17").Value = "FireFlink_26270"
$ws2.Range("M17").Value = "FireFlink_77518"
$ws2.Range("N17").Value = "FireFlink_37757"
$ws2.Range("O17").Value = "FireFlink_87140"
$ws2.Range("P17").Value = "FireFlink_14261"
$ws2.Range("Q17").Value = "FireFlink_51456"
$ws2.Range("R17").Value = "FireFlink_63540"
$ws2.Range("S17").Value = "FireFlink_98531"
$ws2.Range("T17").Value = "FireFlink_19388"
$ws2.Range("A18").Value = "FireFlink_94816"
$ws2.Range("B18").Value = "FireFlink_09536"
$ws2.Range("C18").Value = "FireFlink_47450"
$ws2.Range("D18").Value = "FireFlink_28529"
$ws2.Range("E18").Value = "FireFlink_89297"
$ws2.Range("F18").Value = "FireFlink_47417"
$ws2.Range("G18").Value = "FireFlink_06795"
$ws2.Range("H18").Value = "FireFlink_69692"
$ws2.Range("I18").Value = "FireFlink_18200"
$ws2.Range("J18").Value = "FireFlink_42520"
$ws2.Range("K18").Value = "FireFlink_02356"
$ws2.Range("L18").Value = "FireFlink_73328"
$ws2.Range("M18").Value = "FireFlink_08827"
$ws2.Range("N18").Value = "FireFlink_08556"
$ws2.Range("O18").Value = "FireFlink_81525"
$ws2.Range("P18").Value = "FireFlink_26503"
$ws2.Range("Q18").Value = "FireFlink_85254"
$ws2.Range("R18").Value = "FireFlink_92121"
$ws2.Range("S18").Value = "FireFlink_71659"
$ws2.Range("T18").Value = "FireFlink_47510"
$ws2.Range("A19").Value = "FireFlink_43004"
$ws2.Range("B19").Value = "FireFlink_51746"
$ws2.Range("C19").Value = "FireFlink_14619"
$ws2.Range("D19").Value = "FireFlink_27408"
$ws2.Range("E19").Value = "FireFlink_76858"
$ws2.Range("F19").Value = "FireFlink_00969"
$ws2.Range("G19").Value = "FireFlink_11950"
$ws2.Range("H19").Value = "FireFlink_38581"
$ws2.Range("I19").Value = "FireFlink_64919"
$ws2.Range("J19").Value = "FireFlink_82099"
$ws2.Range("K19").Value = "FireFlink_77926"
$ws2.Range("L19").Value = "FireFlink_28648"
$ws2.Range("M19").Value = "FireFlink_78369"
$ws2.Range("N19").Value = "FireFlink_13397"
$ws2.Range("O19").Value = "FireFlink_55678"
$ws2.Range("P19").Value = "FireFlink_25739"
$ws2.Range("Q19").Value = "FireFlink_33148"
$ws2.Range("R19").Value = "FireFlink_45717"
$ws2.Range("S19").Value = "FireFlink_55450"
$ws2.Range("T19").Value = "FireFlink_94099"
$ws2.Range("A20").Value = "FireFlink_83078"
$ws2.Range("B20").Value = "FireFlink_18163"
$ws2.Range("C20").Value = "FireFlink_16886"
$ws2.Range("D20").Value = "FireFlink_34161"
$ws2.Range("E20").Value = "FireFlink_56748"
$ws2.Range("F20").Value = "FireFlink_25609"
$ws2.Range("G20").Value = "FireFlink_61762"
$ws2.Range("H20").Value = "FireFlink_31484"
$ws2.Range("I20").Value = "FireFlink_03539"
$ws2.Range("J20").Value = "FireFlink_50146"
$ws2.Range("K20").Value = "FireFlink_89262"
$ws2.Range("L20").Value = "FireFlink_44993"
$ws2.Range("M20").Value = "FireFlink_97848"
$ws2.Range("N20").Value = "FireFlink_25224"
$ws2.Range("O20").Value = "FireFlink_43574"
$ws2.Range("P20").Value = "FireFlink_58076"
$ws2.Range("Q20").Value = "FireFlink_12276"
$ws2.Range("R20").Value = "FireFlink_41085"
$ws2.Range("S20").Value = "FireFlink_18875"
$ws2.Range("T20").Value = "FireFlink_73403"
$ws2.Range("A21").Value = "FireFlink_35430"
$ws2.Range("B21").Value = "FireFlink_32872"
$ws2.Range("C21").Value = "FireFlink_54469"
$ws2.Range("D21").Value = "FireFlink_33321"
$ws2.Range("E21").Value = "FireFlink_06302"
$ws2.Range("F21").Value = "FireFlink_96908"
$ws2.Range("G21").Value = "FireFlink_51003"
$ws2.Range("H21").Value = "FireFlink_51092"
$ws2.Range("I21").Value = "FireFlink_16660"
$ws2.Range("J21").Value = "FireFlink_07148"
$ws2.Range("K21").Value = "FireFlink_00777"
$ws2.Range("L21").Value = "FireFlink_40511"
$ws2.Range("M21").Value = "FireFlink_43452"
$ws2.Range("N21").Value = "FireFlink_74147"
$ws2.Range("O21").Value = "FireFlink_22246"
$ws2.Range("P21").Value = "FireFlink_67486"
$ws2.Range("Q21").Value = "FireFlink_16270"
$ws2.Range("R21").Value = "FireFlink_58585"
$ws2.Range("S21").Value = "FireFlink_48162"
$ws2.Range("T21").Value = "FireFlink_58141"
$ws2.Range("A22").Value = "FireFlink_23730"
$ws2.Range("B22").Value = "FireFlink_92862"
$ws2.Range("C22").Value = "FireFlink_76379"
$ws2.Range("D22").Value = "FireFlink_05725"
$ws2.Range("E22").Value = "FireFlink_87305"
$ws2.Range("F22").Value = "FireFlink_17912"
$ws2.Range("G22").Value = "FireFlink_02936"
$ws2.Range("H22").Value = "FireFlink_40580"
$ws2.Range("I22").Value = "FireFlink_42154"
$ws2.Range("J22").Value = "FireFlink_57740"
$ws2.Range("K22").Value = "FireFlink_84146"
$ws2.Range("L22").Value = "FireFlink_03144"
$ws2.Range("M22").Value = "FireFlink_45805"
$ws2.Range("N22").Value = "FireFlink_33287"
$ws2.Range("O22").Value = "FireFlink_86568"
$ws2.Range("P22").Value = "FireFlink_23329"
$ws2.Range("Q22").Value = "FireFlink_88979"
$ws2.Range("R22").Value = "FireFlink_35639"
$ws2.Range("S22").Value = "FireFlink_00073"
$ws2.Range("T22").Value = "FireFlink_71766"
$ws2.Range("A23").Value = "FireFlink_64080"
$ws2.Range("B23").Value = "FireFlink_93137"
$ws2.Range("C23").Value = "FireFlink_50723"
$ws2.Range("D23").Value = "FireFlink_25922"
$ws2.Range("E23").Value = "FireFlink_29745"
$ws2.Range("F23").Value = "FireFlink_79172"
$ws2.Range("G23").Value = "FireFlink_79905"
$ws2.Range("H23").Value = "FireFlink_67811"
$ws2.Range("I23").Value = "FireFlink_74007"
$ws2.Range("J23").Value = "FireFlink_30136"

# Latest scraped project name goes on the visible summary row.
$ws1.Range("B2").Value = "FireFlink_30136"
